# Generate Report for Handback
# Adds a new handback-status row for file 6525353b-4c16-4cfb-be57-d26d14a0f9a4
# to the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$guid = "6525353b-4c16-4cfb-be57-d26d14a0f9a4"
$hash = "773a0f640f5cb27ba8d0c2b3c34d0d2eca32e14e"
$mdName = "$guid.md"

$hlColor = 15570276  # BGR for FF6495ED, matches the workbook's existing HyperLink style
$underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle

function Style-AsHyperlink($range) {
    $range.Font.Underline = $underline
    $range.Font.Color = $hlColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName",
    $null, $null, $mdName) | Out-Null
Style-AsHyperlink $wsOverview.Range("A4")

$wsOverview.Range("B4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C4").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$guid.$hash.zh-cn.xlf"

$wsZh.Range("A4").Value = $mdName
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName",
    $null, $null, $mdName) | Out-Null
Style-AsHyperlink $wsZh.Range("A4")

$wsZh.Range("B4").Value = "Handed back: in sync with en-US"

$wsZh.Range("C4").Value = $zhXlf
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    $null, $null, $zhXlf) | Out-Null
Style-AsHyperlink $wsZh.Range("C4")

$wsZh.Range("D4").Value = "2016-02-29 03:49:06"
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("E4").Value = $mdName
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$mdName",
    $null, $null, $mdName) | Out-Null
Style-AsHyperlink $wsZh.Range("E4")

$wsZh.Range("F4").Value = $zhXlf
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    $null, $null, $zhXlf) | Out-Null
Style-AsHyperlink $wsZh.Range("F4")

$wsZh.Range("G4").Value = "2016-02-29 03:49:50"
$wsZh.Range("H4").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf = "$guid.$hash.de-de.xlf"

$wsDe.Range("A4").Value = $mdName
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName",
    $null, $null, $mdName) | Out-Null
Style-AsHyperlink $wsDe.Range("A4")

$wsDe.Range("B4").Value = "Handed back: in sync with en-US"

$wsDe.Range("C4").Value = $deXlf
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    $null, $null, $deXlf) | Out-Null
Style-AsHyperlink $wsDe.Range("C4")

$wsDe.Range("D4").Value = "2016-02-29 03:49:18"
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("E4").Value = $mdName
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$mdName",
    $null, $null, $mdName) | Out-Null
Style-AsHyperlink $wsDe.Range("E4")

$wsDe.Range("F4").Value = $deXlf
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    $null, $null, $deXlf) | Out-Null
Style-AsHyperlink $wsDe.Range("F4")

$wsDe.Range("G4").Value = "2016-02-29 03:50:11"
$wsDe.Range("H4").Value = "Include"

Write-Host "Row 4 added to Overview, zh-cn and de-de sheets."
